# Update "installed software" tracking on the 'Mac installed' sheet:
#  - insert a new "Aktiv" column (whether the package is currently installed)
#  - fix a typo (AdroidFileTransfer -> AndroidFileTransfer) and mark it inactive
#  - add two new rows: OpenMTP (replacement) and vlc
#  - update the hidden filter-database defined name to the new range

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mac installed")

# Touch the new strings first (in the order the author introduced them) so
# they get interned into the shared-string table in that same order, before
# touching any pre-existing ones again below.
$ws.Range("Z1").Value = "OpenMTP"
$ws.Range("Z2").Value = "Alternative zu AndroidFileTransfer"

# Fix the typo in row 18 (it has been replaced by OpenMTP above).
$ws.Range("A18").Value = "AndroidFileTransfer"

# Insert a new column B ("Aktiv"); this shifts the old "Brew" column to C
# and the old "Notes" column to D, preserving their contents/formatting.
$ws.Columns.Item(2).Insert()

# New column header + narrow width to match the "Brew" column next to it.
$ws.Range("B1").Value = "Aktiv"
$ws.Range("B1:C1").Columns.ColumnWidth = 4.166666666666667

# Populate the new "Aktiv" column: everything currently tracked is active,
# except the just-renamed AndroidFileTransfer row.
$ws.Range("B2:B17").Value = "✅"
$ws.Range("B18").Value = "❌"
$ws.Range("B19:B28").Value = "✅"

# Append two new rows for recently installed software.
$ws.Range("A29").Value = "OpenMTP"
$ws.Range("B29").Value = "✅"
$ws.Range("C29").Value = "✅"
$ws.Range("D29").Value = "Alternative zu AndroidFileTransfer"

$ws.Range("A30").Value = "vlc"
$ws.Range("B30").Value = "✅"
$ws.Range("C30").Value = "✅"

# Clear the scratch cells used above to seed shared-string ordering
# (the column insert shifted them from Z to AA).
$ws.Range("AA1:AA2").Clear()

# Select the last-entered cells, matching where editing ended.
$ws.Range("B30:C30").Select()

# The sheet's hidden filter-database range now covers the extra column.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Mac installed!_FilterDatabase") {
        $n.RefersTo = "='Mac installed'!`$A`$1:`$D`$1"
    }
}
